# Update "想去人数" (number of people interested) counts on the
# "展览" (Exhibition) and "全部类型" (All types) sheets.
# These two sheets list the same events (全部类型 just interleaves the
# 演出/Performance rows among them), so the same F-column increases are
# applied to both, at their respective row numbers.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet1.Range("F2").Value = 12524
$sheet1.Range("F3").Value = 599
$sheet1.Range("F4").Value = 2034
$sheet1.Range("F6").Value = 385
$sheet1.Range("F8").Value = 12493
$sheet1.Range("F9").Value = 3077
$sheet1.Range("F17").Value = 2836
$sheet1.Range("F18").Value = 6079

$sheet4 = $wb.Worksheets.Item("全部类型")
$sheet4.Range("F2").Value = 12524
$sheet4.Range("F3").Value = 599
$sheet4.Range("F4").Value = 2034
$sheet4.Range("F7").Value = 385
$sheet4.Range("F9").Value = 12493
$sheet4.Range("F10").Value = 3077
$sheet4.Range("F18").Value = 2836
$sheet4.Range("F20").Value = 6079
